# Automatische test-sync: 2025-07-23 22:24:50
# Append the 5th test-mail entry to the "Logs" sheet and the matching
# rollup row to the "Dashboard" sheet, then extend the dependent
# conditional-formatting ranges and the bar-chart series references.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs!A15:J15 - new row of data
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A15").Value = "Ik probeer in te loggen maar krijg steeds een foutmelding."
$logs.Range("B15").Value = "mailmind.test@zohomail.eu"
$logs.Range("C15").Value = "Testmail #5: Ik probeer in te loggen maar krijg steeds een foutmelding."
$logs.Range("D15").Value = "IT / Technisch probleem"

$antwoord15 = @"
Beste klant,
Bedankt voor het melden van dit probleem. Om u beter van dienst te kunnen zijn, vragen wij u vriendelijk om de foutmelding die u ontvangt door te sturen, samen met uw gebruikersnaam. Zo kunnen wij het probleem nauwkeuriger onderzoeken en u een oplossing bieden.
Met vriendelijke groet,
[Naam bedrijf] E-mailassistent
"@
$logs.Range("E15").Value = $antwoord15

$logs.Range("F15").Value = "2025-07-23 22:24:46"
$logs.Range("G15").Value = "Ja"
$logs.Range("H15").Value = "Nee"
$logs.Range("I15").Value = "Ja"
$logs.Range("J15").Value = "Nee"

# Keep the row at the sheet's default height (matches the other data
# rows, which carry no explicit row height) instead of the taller,
# auto-computed height Excel assigns for the wrapped multi-line text.
$logs.Rows.Item(15).AutoFit()

# ---------------------------------------------------------------------
# 2. Extend the conditional-formatting ranges on Logs from row 14 to 15
# ---------------------------------------------------------------------
$logs.Range("D2:D14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D15"))
$logs.Range("G2:G14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G15"))
$logs.Range("H2:H14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H15"))
$logs.Range("I2:I14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I15"))
$logs.Range("J2:J14").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J15"))

# ---------------------------------------------------------------------
# 3. Dashboard!A6:B6 - new rollup row for the "IT / Technisch probleem"
#    category
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("A6").Value = "IT / Technisch probleem"
$dash.Range("B6").Value = 1

# ---------------------------------------------------------------------
# 4. Extend the chart's category/value series references to include the
#    new Dashboard row 6
# ---------------------------------------------------------------------
$chart = $dash.ChartObjects(1).Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$6,Dashboard!`$B`$2:`$B`$6,1)"
